# Update status values for existing candidates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = "1st Interview"
$ws.Range("E19").Value = "2nd Interview"
$ws.Range("E22").Value = "4th Interview"
$ws.Range("E23").Value = "1st Interview"

# Add a new candidate row
$ws.Range("A27").Value = 871
$ws.Range("B27").Value = "LaunchDarkly"
$ws.Range("C27").Value = "Strategic AE U.S"
$ws.Range("D27").Value = "Paul Plofchan"
$ws.Range("E27").Value = "CV Sent"
